# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Numeric-looking "Price" column cells are forced to Text format before the
# write (then reset to the Normal style) so they keep their exact decimal
# string representation instead of being auto-coerced to a Number by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.538.70"
$ws.Range("D3").Value = "3.146.44"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +11.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.33"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.424"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.06%  "
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").Value = "3.688.58"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000170"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.78%  "
$ws.Range("D16").Value = "58.590.92"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").Value = "3.147.17"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.23%  "
$ws.Range("D29").Value = "0.0₃0855"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.75%  "
$ws.Range("E34").Value = "  +1.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.43%  "
$ws.Range("D39").Value = "2.648.30"
$ws.Range("E39").Value = "  +10.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0682"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.04%  "
$ws.Range("E42").Value = "  +3.99%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.708"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.64%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.21%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0285"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.52%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "3.189.12"
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.102"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +12.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.977"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.33%  "
